$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Total Line of Code" formulas in column C: add the new 350 LOC reflex layer ---
$ws.Range("C4").Formula = "= 136 + 46 + 350"
$ws.Range("C5").Formula = "= 137 + 60 + 350"
$ws.Range("C6").Formula = "= 184 + 3 + 350"
$ws.Range("C7").Formula = "= 137 +  225 + 350"
$ws.Range("C8").Formula = "= 149 + 239 + 350"

# --- Update "Man Hours" values in column E to reflect the additional work ---
$ws.Range("E5").Value = 32
$ws.Range("E7").Value = 25
$ws.Range("E8").Value = 30
$ws.Range("E9").Value = 16

# --- Move the active selection to G21 (matches the saved view state) ---
$ws.Range("G21").Select()
